# Update agencycount_compare.xlsx with new on_multiple_teams field data.
# For each changed agency row, the current_count (B) increases and the
# change (D) column is updated to reflect the new difference from
# previous_count (C), which is B - C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new current_count (B) value.
# The "change" column (D) is recalculated as B - C for each row.
$updates = @{
    2  = 8
    3  = 10
    9  = 21
    10 = 43
    11 = 33
    12 = 21
    14 = 28
    15 = 31
    19 = 15
    20 = 13
    26 = 25
    27 = 19
    35 = 18
}

foreach ($row in $updates.Keys) {
    $newB = $updates[$row]
    $cVal = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 2).Value = $newB
    $ws.Cells.Item($row, 4).Value = $newB - $cVal
}
